$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1213.3182
$ws.Range("J19").Value = 1188.75
$ws.Range("L19").Value = 1188.75
$ws.Range("N19").Value = -1538.75

# Row 29
$ws.Range("H29").Value = 10
$ws.Range("I29").Value = 10
$ws.Range("K29").Value = 30
$ws.Range("M29").Value = 251

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 53
$ws.Range("H53").Value = 233.36363
$ws.Range("I53").Value = 203.14285
$ws.Range("J53").Value = 286.25
$ws.Range("K53").Value = 203.14285
$ws.Range("L53").Value = 286.25
$ws.Range("M53").Value = 433.85715
$ws.Range("N53").Value = -1560.25

# Row 74
$ws.Range("H74").Value = 57156.434
$ws.Range("I74").Value = 85415.30499999999
$ws.Range("K74").Value = 85415.30499999999
$ws.Range("M74").Value = -84479.30499999999

# Row 77
$ws.Range("H77").Value = 57156.434
$ws.Range("I77").Value = 85415.30499999999
$ws.Range("K77").Value = 427076.525
$ws.Range("M77").Value = -422396.525

# Row 99
$ws.Range("H99").Value = 2401.75
$ws.Range("I99").Value = 410.66666
$ws.Range("J99").Value = 4392.8335
$ws.Range("K99").Value = 1231.99998
$ws.Range("L99").Value = 13178.5005
$ws.Range("M99").Value = 266.0000199999999
$ws.Range("N99").Value = -16174.5005

# Row 100
$ws.Range("H100").Value = 1570.5714
$ws.Range("I100").Value = 1570.5714
$ws.Range("K100").Value = 1570.5714
$ws.Range("M100").Value = -1029.5714

# Row 113
$ws.Range("H113").Value = 75108.31
$ws.Range("J113").Value = 17291.5
$ws.Range("L113").Value = 17291.5
$ws.Range("N113").Value = -23799.5

# Row 116
$ws.Range("H116").Value = 9302
$ws.Range("I116").Value = 10002.875
$ws.Range("J116").Value = 6498.5
$ws.Range("K116").Value = 10002.875
$ws.Range("L116").Value = 6498.5
$ws.Range("M116").Value = -6560.875
$ws.Range("N116").Value = -13382.5

# Row 132
$ws.Range("H132").Value = 60401.176
$ws.Range("I132").Value = 38579.258
$ws.Range("K132").Value = 115737.774
$ws.Range("M132").Value = -113207.774

# Row 137
$ws.Range("H137").Value = 2174.682
$ws.Range("I137").Value = 2206.875
$ws.Range("J137").Value = 2088.8333
$ws.Range("K137").Value = 6620.625
$ws.Range("L137").Value = 6266.499899999999
$ws.Range("M137").Value = -4070.625
$ws.Range("N137").Value = -11366.4999

# Row 141
$ws.Range("H141").Value = 1297.5
$ws.Range("I141").Value = 1096.6666
$ws.Range("K141").Value = 3289.9998
$ws.Range("M141").Value = 1890.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2685.9387
$ws.Range("I32").Value = 2162.4666
$ws.Range("K32").Value = 2162.4666
$ws.Range("M32").Value = -1875.4666

# Row 45
$ws.Range("H45").Value = 9406.267
$ws.Range("I45").Value = 11251.454
$ws.Range("J45").Value = 4332
$ws.Range("K45").Value = 11251.454
$ws.Range("L45").Value = 4332
$ws.Range("M45").Value = -10874.454
$ws.Range("N45").Value = -5086

# Row 82
$ws.Range("H82").Value = 39998.75
$ws.Range("J82").Value = 39998.75
$ws.Range("L82").Value = 39998.75
$ws.Range("N82").Value = -40720.75

# Row 85
$ws.Range("H85").Value = 39998.75
$ws.Range("J85").Value = 39998.75
$ws.Range("L85").Value = 39998.75
$ws.Range("N85").Value = -42494.75

# Row 97
$ws.Range("H97").Value = 3701.3333
$ws.Range("I97").Value = 3841.6
$ws.Range("K97").Value = 3841.6
$ws.Range("M97").Value = -3345.6

# Row 102
$ws.Range("H102").Value = 7125
$ws.Range("I102").Value = 9250
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 9250
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -7628
$ws.Range("N102").Value = -8244

# Row 110
$ws.Range("H110").Value = 15673.5
$ws.Range("I110").Value = 26825.143
$ws.Range("K110").Value = 26825.143
$ws.Range("M110").Value = -24780.143

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1811.8948
$ws.Range("I94").Value = 2132.4167
$ws.Range("K94").Value = 2132.4167
$ws.Range("M94").Value = -1681.4167

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 105
$ws.Range("H105").Value = 1680.6923
$ws.Range("I105").Value = 1666.125
$ws.Range("J105").Value = 1704
$ws.Range("K105").Value = 1666.125
$ws.Range("L105").Value = 1704
$ws.Range("M105").Value = 80.875
$ws.Range("N105").Value = -5198

# Row 113
$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# Row 122
$ws.Range("H122").Value = 2175.3333
$ws.Range("I122").Value = 1756
$ws.Range("J122").Value = 3014
$ws.Range("K122").Value = 5268
$ws.Range("L122").Value = 9042
$ws.Range("M122").Value = -2818
$ws.Range("N122").Value = -13942

# Row 134
$ws.Range("H134").Value = 18618.56
$ws.Range("I134").Value = 8798.200000000001
$ws.Range("K134").Value = 26394.6
$ws.Range("M134").Value = -23859.6

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 1397
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1397
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4191
$ws.Range("N12").Value = -4537
$ws.Range("M12").ClearContents()

# Row 102
$ws.Range("H102").Value = 8750

# Row 131
$ws.Range("H131").Value = 30030.25
$ws.Range("J131").Value = 2471.7188
$ws.Range("L131").Value = 7415.1564
$ws.Range("N131").Value = -17495.1564

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 2000000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# Row 34
$ws.Range("H34").Value = 47000
$ws.Range("J34").Value = 47000
$ws.Range("L34").Value = 47000
$ws.Range("N34").Value = -47536

# Row 43
$ws.Range("H43").Value = 14033.4
$ws.Range("I43").Value = 13791.75
$ws.Range("K43").Value = 13791.75
$ws.Range("M43").Value = -13640.75

# Row 76
$ws.Range("H76").Value = 47000
$ws.Range("J76").Value = 47000
$ws.Range("L76").Value = 47000
$ws.Range("N76").Value = -47630

# Row 79
$ws.Range("H79").Value = 47000
$ws.Range("J79").Value = 47000
$ws.Range("L79").Value = 47000
$ws.Range("N79").Value = -49184

# Row 107
$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 5000
$ws.Range("K107").Value = 5000
$ws.Range("M107").Value = -3080

# Row 113
$ws.Range("H113").Value = 1955.5
$ws.Range("I113").Value = 1955.5
$ws.Range("K113").Value = 1955.5
$ws.Range("M113").Value = 214.5

# Row 122
$ws.Range("H122").Value = 899
$ws.Range("I122").Value = 899
$ws.Range("K122").Value = 2697
$ws.Range("M122").Value = -247

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 80000
$ws.Range("J3").Value = 80000
$ws.Range("L3").Value = 80000
$ws.Range("N3").Value = -80224

# Row 7
$ws.Range("H7").Value = 11945.728
$ws.Range("I7").Value = 13425.5
$ws.Range("J7").Value = 7999.6665
$ws.Range("K7").Value = 13425.5
$ws.Range("L7").Value = 7999.6665
$ws.Range("M7").Value = -13313.5
$ws.Range("N7").Value = -8223.666499999999

# Row 15
$ws.Range("H15").Value = 80000
$ws.Range("J15").Value = 80000
$ws.Range("L15").Value = 80000
$ws.Range("N15").Value = -80340

# Row 22
$ws.Range("H22").Value = 1547.5
$ws.Range("J22").Value = 1867.3334
$ws.Range("L22").Value = 1867.3334
$ws.Range("N22").Value = -2457.3334

# Row 27
$ws.Range("H27").Value = 1547.5
$ws.Range("J27").Value = 1867.3334
$ws.Range("L27").Value = 1867.3334
$ws.Range("N27").Value = -2081.3334

# Row 40
$ws.Range("H40").Value = 3355
$ws.Range("I40").Value = 3096.9443
$ws.Range("K40").Value = 3096.9443
$ws.Range("M40").Value = -2960.9443

# Row 100
$ws.Range("H100").Value = 29666.666
$ws.Range("I100").Value = 27500
$ws.Range("J100").Value = 34000
$ws.Range("K100").Value = 27500
$ws.Range("L100").Value = 34000
$ws.Range("M100").Value = -26959
$ws.Range("N100").Value = -35082

# Row 122
$ws.Range("H122").Value = 4779.067
$ws.Range("I122").Value = 4140.5
$ws.Range("K122").Value = 12421.5
$ws.Range("M122").Value = -9971.5

# Row 126
$ws.Range("H126").Value = 11945.728
$ws.Range("I126").Value = 13425.5
$ws.Range("J126").Value = 7999.6665
$ws.Range("K126").Value = 40276.5
$ws.Range("L126").Value = 23998.9995
$ws.Range("M126").Value = -37806.5
$ws.Range("N126").Value = -28938.9995

# Row 132
$ws.Range("H132").Value = 6083.1113
$ws.Range("I132").Value = 5345.846
$ws.Range("K132").Value = 16037.538
$ws.Range("M132").Value = -13507.538

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 22612.25
$ws.Range("I75").Value = 18999.5
$ws.Range("J75").Value = 26225
$ws.Range("K75").Value = 18999.5
$ws.Range("L75").Value = 26225
$ws.Range("M75").Value = -18063.5
$ws.Range("N75").Value = -28097

# Row 78
$ws.Range("H78").Value = 22612.25
$ws.Range("I78").Value = 18999.5
$ws.Range("J78").Value = 26225
$ws.Range("K78").Value = 56998.5
$ws.Range("L78").Value = 78675
$ws.Range("M78").Value = -52318.5
$ws.Range("N78").Value = -88035

# Row 107
$ws.Range("H107").Value = 519.2
$ws.Range("I107").Value = 533
$ws.Range("K107").Value = 1599
$ws.Range("M107").Value = 321

# Row 126
$ws.Range("H126").Value = 2644.8076
$ws.Range("I126").Value = 2643.3333
$ws.Range("J126").Value = 2648.125
$ws.Range("K126").Value = 7929.999899999999
$ws.Range("L126").Value = 7944.375
$ws.Range("M126").Value = -5459.999899999999
$ws.Range("N126").Value = -12884.375

# Row 132
$ws.Range("H132").Value = 3693.524
$ws.Range("I132").Value = 3426.389
$ws.Range("K132").Value = 10279.167
$ws.Range("M132").Value = -7749.167000000001
